$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 1 (header row): extend the header style (bold font + thin box
# border, same as the existing B1:G1 headers) across the new H1:N1
# columns before writing their text.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ------------------------------------------------------------------
# Row 2 (data row): extend the data-row style across the new H2:N2
# columns, same as it is used by the pre-existing C2:G2 cells.
# ------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("C2").Value = "蘇震清"
$ws.Range("D2").Value = "臺灣新光商業銀行五常分行臺北市中山區龍江路356巷76號"
$ws.Range("E2").Value = 1468576
$ws.Range("F2").Value = "100年05月16日"
$ws.Range("G2").Value = "貸款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"

# J2 holds the literal text "2012-04-30" (not a date serial), so force
# a text number format before assigning it, then restore the normal
# data-row formatting on top (keeps the cell out of "Date" display).
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("C2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$ws.Range("K2").Value = "蘇震清"
$ws.Range("L2").Value = 1718
$ws.Range("M2").Value = "tmp31e11"
$ws.Range("N2").Value = 13
